$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add I0 / IF headers (copy the header style/format from H1, same as the other headers)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:J63 data values (row -> [I, J])
$ijData = @(
    @(1,2),
    @(6,7),
    @(7,7),
    @(7,8),
    @(1,2),
    @(7,8),
    @(5,5),
    @(5,5),
    @(9,9),
    @(7,7),
    @(7,7),
    @(1,2),
    @(8,8),
    @(6,7),
    @(6,6),
    @(5,6),
    @(8,8),
    @(6,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(9,9),
    @(7,7),
    @(7,7),
    @(9,9),
    @(6,7),
    @(9,9),
    @(6,6),
    @(10,10),
    @(3,4),
    @(5,6),
    @(6,6),
    @(7,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(3,4),
    @(6,7),
    @(8,8),
    @(10,11),
    @(7,7),
    @(9,9),
    @(9,9),
    @(7,8),
    @(6,6),
    @(4,5),
    @(8,9),
    @(5,6),
    @(5,6),
    @(9,9),
    @(6,7),
    @(4,4),
    @(4,4),
    @(8,8),
    @(8,8),
    @(7,7),
    @(7,7)
)

for ($k = 0; $k -lt $ijData.Length; $k++) {
    $row = $k + 2
    $pair = $ijData[$k]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
